$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format for Price/Volume columns so that
# numeric-looking strings (e.g. "11.00", "0.0900") keep their exact text
# representation instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.181.16"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "2.306.27"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "312.85"
$ws.Range("E5").Value = "  -3.35%  "
$ws.Range("D6").Value = "105.71"
$ws.Range("E6").Value = "  +5.73%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").Value = "40.28"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "0.977"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").Value = "15.56"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("D16").Value = "2.655.76"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").Value = "2.309.01"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "42.136.57"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").Value = "7.68"
$ws.Range("E19").Value = "  -4.52%  "
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "74.56"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("E22").Value = "  -5.76%  "
$ws.Range("D23").Value = "259.93"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "2.30"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "9.34"
$ws.Range("E25").Value = "  -6.86%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "11.00"
$ws.Range("E27").Value = "  -3.57%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("D29").Value = "22.75"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").Value = "35.83"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").Value = "0.0900"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "163.76"
$ws.Range("E32").Value = "  -6.50%  "
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  -5.39%  "
$ws.Range("D34").Value = "5.85"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  +12.68%  "
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").Value = "0.0352"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("D41").Value = "71.93"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").Value = "98.54"
$ws.Range("E42").Value = "  +8.33%  "
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").Value = "0.228"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "12.29"
$ws.Range("E46").Value = "  +4.70%  "
$ws.Range("D47").Value = "112.52"
$ws.Range("E47").Value = "  -4.74%  "
$ws.Range("D48").Value = "9.03"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("D50").Value = "74.52"
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("E51").Value = "  +0.12%  "
